# "complete add hero without command pattern"
# Checklist sheet: column B marks an item "OK" once implemented.
# Row 4  -> "2. Add hero to player"
# Row 6  -> "4. display all player"
# Both are now complete, so mark them "OK".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = "OK"
$ws.Range("B6").Value = "OK"

# Leave the selection where the author ended up editing.
$null = $ws.Range("B14").Select()
